$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (Sandía @ Vega Monumental Concepción) needs to be
# inserted chronologically between the existing row 32 and row 33 blocks.
# Insert 3 blank rows at 33:35 (this shifts the old rows 33:67 down to 36:70
# and keeps their original values, styles, dates, etc. intact).
$ws.Rows("33:35").Insert()

# Populate the 3 newly inserted rows with the new observations.
# Row 33: Extra
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44553
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112028
$ws.Range("G33").Value = "Sandia"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Extra"
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 3400
$ws.Range("L33").Value = 3400
$ws.Range("M33").Value = 3400
$ws.Range("N33").Value = "`$/unidad"
$ws.Range("O33").Value = "Región de O'Higgins"
$ws.Range("P33").Value = 3400
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"

# Row 34: Primera
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44553
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112028
$ws.Range("G34").Value = "Sandia"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 500
$ws.Range("K34").Value = 2800
$ws.Range("L34").Value = 2800
$ws.Range("M34").Value = 2800
$ws.Range("N34").Value = "`$/unidad"
$ws.Range("O34").Value = "Región de O'Higgins"
$ws.Range("P34").Value = 2800
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"

# Row 35: Segunda
$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value = 44553
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = 100112028
$ws.Range("G35").Value = "Sandia"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 2400
$ws.Range("L35").Value = 2400
$ws.Range("M35").Value = 2400
$ws.Range("N35").Value = "`$/unidad"
$ws.Range("O35").Value = "Región de O'Higgins"
$ws.Range("P35").Value = 2400
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"
